$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grade values in row 10
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5
$ws.Range("M10").Value = 5

# Update grade values in row 21
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 5

# Update frozen pane top-left cell and selection to reflect new scroll position
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("J10").Select()
